# Update the "built on" timestamp embedded in the version string wherever it
# appears in the workbook (About sheet notes + per-row build_version column
# on the Boundaries and methane sources sheet).

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($val -ne $null -and $val -is [string] -and $val.Contains($oldStamp)) {
                $cell.Value = $val.Replace($oldStamp, $newStamp)
            }
        }
    }
}
